$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $value) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("A2").Value = $true
$ws.Range("B2").Value = "edit"
Set-TextCell 2 "2370"

$ws.Range("A3").Value = $false
$ws.Range("B3").Value = "edit"
Set-TextCell 3 "2353"

$ws.Range("A4").Value = $false
$ws.Range("B4").Value = "categorize"
Set-TextCell 4 "1613"

$ws.Range("A5").Value = $true
$ws.Range("B5").Value = "categorize"
Set-TextCell 5 "1278"

$ws.Range("A6").Value = $true
$ws.Range("B6").Value = "new"
Set-TextCell 6 "714"

$ws.Range("A7").Value = $false
$ws.Range("B7").Value = "log"
Set-TextCell 7 "366"

$ws.Range("A8").Value = $true
$ws.Range("B8").Value = "log"
Set-TextCell 8 "175"

$ws.Range("A9").Value = $false
$ws.Range("B9").Value = "new"
Set-TextCell 9 "129"
